$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeaponData")

# --- New column headers (English), row 1 ---
$ws.Range("E1").Value = "UIPath"
$ws.Range("F1").Value = "Description"
# Match the header formatting used by the existing header cells (A1:D1)
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# --- New column headers (Chinese), row 2 ---
$ws.Range("E2").Value = "武器UI路徑"
$ws.Range("F2").Value = "武器說明"

# --- New column type markers, row 3 ---
$ws.Range("E3").Value = "string"
$ws.Range("F3").Value = "string"

# --- New data values for the MagicBall weapon row, row 4 ---
$ws.Range("E4").Value = "Assets/ArtResources/Weapons/Weapons Sprite Sheet.png[Weapons Sprite Sheet_123]"
$ws.Range("F4").Value = "朝向の方向に素早く発射します"

# --- Column widths for the two new columns ---
$ws.Range("E:E").ColumnWidth = 66.16666666666667
$ws.Range("F:F").ColumnWidth = 25.833333333333332
